$wb = $excel.ActiveWorkbook

# --- Functional Requirements sheet: add row 13 (report-issue requirement) ---
$wsFunc = $wb.Worksheets.Item("Functional Requirements")
$wsFunc.Range("B17").Value = 13
$wsFunc.Range("C17").Value = "The TVM should have a report issue on the menu"
$wsFunc.Range("D17").Value = "Want"
$wsFunc.Range("E17").Value = "The TVM will have a report issue option for user to send the issue to the admin"

# --- Non-Functional Requirements sheet: normalize row 11 formatting to match
#     the rest of the table (drop the stray "applyProtection" style variant by
#     copying the format from the row above) ---
$wsNonFunc = $wb.Worksheets.Item("Non-Functional Requirements")
$wsNonFunc.Range("B10:E10").Copy()
$wsNonFunc.Range("B11:E11").PasteSpecial(-4122)
